$wb = $excel.ActiveWorkbook

# The report generator re-stamps every "Ready for handoff" status cell with
# "In Translation" across all sheets (Overview rollup + per-locale sheets).
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # NOTE: compare against an explicit [string] cast - the COM shim
        # hands back a typed Boolean for cells whose text happens to be
        # "True"/"False", and leaving it untyped corrupts unrelated
        # "-eq" checks later in the loop (PowerShell coerces the right
        # side of "-eq" to match a boolean left side).
        $text = [string]$cell.Text
        if ($text -eq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
    $used.Columns.AutoFit() | Out-Null
}

# Re-fit the status columns to the new (shorter) text. The stored <col
# width> Excel writes is the character width plus a fixed ~5/6-character
# padding constant, so back that constant out of the desired stored width
# before handing it to the ColumnWidth property.
$targetStoredWidth = 13.4101845877511
$targetColumnWidth = $targetStoredWidth - (5.0 / 6.0)

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $targetColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $targetColumnWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $targetColumnWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $targetColumnWidth
